$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 36 (existing rows 36-48 shift down to 38-50)
$ws.Rows.Item(36).Resize(2).Insert()

# ---- Set values first, in the same order the author typed them, so the
# ---- shared-string table is appended in the matching order (A36, A37, B37, B36) ----
$a36 = $ws.Cells.Item(36,1)
$a36.Value = "transform: scale3D(1.2, 1.5, 2)"

$a37 = $ws.Cells.Item(37,1)
$a37.Value = "transform: scaleZ(2)"

$b37 = $ws.Cells.Item(37,2)
$b37.Value = "هاي نفس الفوكاها بس نختصر نكتب بس زي معناته البعد الثالث"

$b36 = $ws.Cells.Item(36,2)
$b36.Value = "هيج نسوي سكيل بالبعد الثالث بس البعد الثالث ميبين الا اكو وياه  فد نوع `nترانسفورم ثاني مثل روتيت او ترانسلايت"

# ---- Row 36 formatting: transform: scale3D(1.2, 1.5, 2) ----
$a36.Font.Name = "Consolas"
$a36.Font.Size = 18
$a36.Font.Bold = $true
$a36.Font.ThemeColor = 1
$a36.Interior.Pattern = -4142
$a36.Interior.ThemeColor = 7
$a36.Interior.TintAndShade = 0.59999389629810485
$a36.VerticalAlignment = -4108
$a36.Borders.Item(9).LineStyle = 1
$a36.Borders.Item(9).Weight = 2
$a36.Borders.Item(9).ColorIndex = -4105

$b36.Font.Name = "Arial"
$b36.Font.Size = 18
$b36.Font.Bold = $false
$b36.HorizontalAlignment = -4108
$b36.VerticalAlignment = -4108
$b36.WrapText = $true

$ws.Rows.Item(36).RowHeight = 61.8

# ---- Row 37 formatting: transform: scaleZ(2) ----
$a37.Font.Name = "Consolas"
$a37.Font.Size = 18
$a37.Font.Bold = $true
$a37.Font.ThemeColor = 1
$a37.Interior.Pattern = -4142
$a37.Interior.ThemeColor = 7
$a37.Interior.TintAndShade = 0.59999389629810485
$a37.HorizontalAlignment = -4131
$a37.VerticalAlignment = -4108

$b37.Font.Name = "Arial"
$b37.Font.Size = 18
$b37.Font.Bold = $false
$b37.HorizontalAlignment = -4108
$b37.VerticalAlignment = -4108
$b37.WrapText = $false

$ws.Rows.Item(37).RowHeight = 64.2

# ---- Update view state to match the committed selection/scroll position ----
$ws.Range("A38").Select()
$excel.ActiveWindow.ScrollRow = 32
